$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data to append as rows 606-613 (columns A-K)
$data = @(
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRC08418","Hemmeter","DRIVING IN MARKED LANES","4511.33","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRC08418","Hemmeter","TURN AND STOP SIGNALS","No Data","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRC08418","Hemmeter","OVI ALCOHOL / DRUGS 1ST","4511.19A1A*","M1","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01291","Hemmeter","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01291","Bunner","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01291","Bunner","PERMISSION REQ'D TO USE LICENSED DOCK","1501:46-12-04","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21CRB01268","Hemmeter","POSSESSION DRUG PARAPHERNALIA","2925.14(C)","M4","No Contest","Guilty","$ 0","$ 0","None","None")
)

# Columns whose literal text values would otherwise be auto-converted by Excel
# into a number or currency value (e.g. "4511.33" -> 4511.33, "$ 0" -> 0).
# For these cells only, temporarily mark the cell as Text before assigning the
# value, then restore the default "Normal" style so no stray cell formatting
# is left behind (matching the plain, unstyled cells used elsewhere in the sheet).
$textColumns = @(4, 8, 9)

$startRow = 606
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        $value = $rowData[$c]
        if ($textColumns -contains $col -and $value -match '^\$?\s?-?\d+(\.\d+)?$') {
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
